$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("BD_Times")
$ws2 = $wb.Worksheets.Item("BD_Jogo")

# --- BD_Times: append rows 20-37 (new round of team stats) ---
$ws1.Cells.Item(20,1).Value = "Metz"
$ws1.Cells.Item(20,2).Value = 1
$ws1.Cells.Item(20,3).Value = 1
$ws1.Cells.Item(20,4).Value = 1
$ws1.Cells.Item(20,5).Value = 1
$ws1.Cells.Item(20,6).Value = 2
$ws1.Cells.Item(20,7).Value = 2
$ws1.Cells.Item(20,8).Value = 4
$ws1.Cells.Item(20,9).Value = 9

$ws1.Cells.Item(21,1).Value = "Marseille"
$ws1.Cells.Item(21,2).Value = 0
$ws1.Cells.Item(21,3).Value = 1
$ws1.Cells.Item(21,4).Value = 1
$ws1.Cells.Item(21,5).Value = 1
$ws1.Cells.Item(21,6).Value = 2
$ws1.Cells.Item(21,7).Value = 2
$ws1.Cells.Item(21,8).Value = 9
$ws1.Cells.Item(21,9).Value = 4

$ws1.Cells.Item(22,1).Value = "Lyon"
$ws1.Cells.Item(22,2).Value = 1
$ws1.Cells.Item(22,3).Value = 1
$ws1.Cells.Item(22,4).Value = 1
$ws1.Cells.Item(22,5).Value = 1
$ws1.Cells.Item(22,6).Value = 1
$ws1.Cells.Item(22,7).Value = 4
$ws1.Cells.Item(22,8).Value = 6
$ws1.Cells.Item(22,9).Value = 2

$ws1.Cells.Item(23,1).Value = "Montpelier"
$ws1.Cells.Item(23,2).Value = 0
$ws1.Cells.Item(23,3).Value = 1
$ws1.Cells.Item(23,4).Value = 1
$ws1.Cells.Item(23,5).Value = 1
$ws1.Cells.Item(23,6).Value = 4
$ws1.Cells.Item(23,7).Value = 1
$ws1.Cells.Item(23,8).Value = 2
$ws1.Cells.Item(23,9).Value = 6

$ws1.Cells.Item(24,1).Value = "Toulouse"
$ws1.Cells.Item(24,2).Value = 1
$ws1.Cells.Item(24,3).Value = 1
$ws1.Cells.Item(24,4).Value = 1
$ws1.Cells.Item(24,5).Value = 1
$ws1.Cells.Item(24,6).Value = 1
$ws1.Cells.Item(24,7).Value = 1
$ws1.Cells.Item(24,8).Value = 5
$ws1.Cells.Item(24,9).Value = 5

$ws1.Cells.Item(25,1).Value = "PSG"
$ws1.Cells.Item(25,2).Value = 0
$ws1.Cells.Item(25,3).Value = 1
$ws1.Cells.Item(25,4).Value = 1
$ws1.Cells.Item(25,5).Value = 1
$ws1.Cells.Item(25,6).Value = 1
$ws1.Cells.Item(25,7).Value = 1
$ws1.Cells.Item(25,8).Value = 5
$ws1.Cells.Item(25,9).Value = 5

$ws1.Cells.Item(26,1).Value = "Lille"
$ws1.Cells.Item(26,2).Value = 1
$ws1.Cells.Item(26,3).Value = 1
$ws1.Cells.Item(26,4).Value = 0
$ws1.Cells.Item(26,5).Value = 0
$ws1.Cells.Item(26,6).Value = 2
$ws1.Cells.Item(26,7).Value = 0
$ws1.Cells.Item(26,8).Value = 7
$ws1.Cells.Item(26,9).Value = 5

$ws1.Cells.Item(27,1).Value = "Nantes"
$ws1.Cells.Item(27,2).Value = 0
$ws1.Cells.Item(27,3).Value = 0
$ws1.Cells.Item(27,4).Value = 1
$ws1.Cells.Item(27,5).Value = 0
$ws1.Cells.Item(27,6).Value = 0
$ws1.Cells.Item(27,7).Value = 2
$ws1.Cells.Item(27,8).Value = 5
$ws1.Cells.Item(27,9).Value = 7

$ws1.Cells.Item(28,1).Value = "Le Havre"
$ws1.Cells.Item(28,2).Value = 1
$ws1.Cells.Item(28,3).Value = 1
$ws1.Cells.Item(28,4).Value = 1
$ws1.Cells.Item(28,5).Value = 1
$ws1.Cells.Item(28,6).Value = 1
$ws1.Cells.Item(28,7).Value = 2
$ws1.Cells.Item(28,8).Value = 5
$ws1.Cells.Item(28,9).Value = 3

$ws1.Cells.Item(29,1).Value = "Brest"
$ws1.Cells.Item(29,2).Value = 0
$ws1.Cells.Item(29,3).Value = 1
$ws1.Cells.Item(29,4).Value = 1
$ws1.Cells.Item(29,5).Value = 1
$ws1.Cells.Item(29,6).Value = 2
$ws1.Cells.Item(29,7).Value = 1
$ws1.Cells.Item(29,8).Value = 3
$ws1.Cells.Item(29,9).Value = 5

$ws1.Cells.Item(30,1).Value = "Lorient"
$ws1.Cells.Item(30,2).Value = 1
$ws1.Cells.Item(30,3).Value = 1
$ws1.Cells.Item(30,4).Value = 1
$ws1.Cells.Item(30,5).Value = 1
$ws1.Cells.Item(30,6).Value = 1
$ws1.Cells.Item(30,7).Value = 1
$ws1.Cells.Item(30,8).Value = 0
$ws1.Cells.Item(30,9).Value = 9

$ws1.Cells.Item(31,1).Value = "Nice"
$ws1.Cells.Item(31,2).Value = 0
$ws1.Cells.Item(31,3).Value = 1
$ws1.Cells.Item(31,4).Value = 1
$ws1.Cells.Item(31,5).Value = 1
$ws1.Cells.Item(31,6).Value = 1
$ws1.Cells.Item(31,7).Value = 1
$ws1.Cells.Item(31,8).Value = 9
$ws1.Cells.Item(31,9).Value = 0

$ws1.Cells.Item(32,1).Value = "Reims"
$ws1.Cells.Item(32,2).Value = 1
$ws1.Cells.Item(32,3).Value = 1
$ws1.Cells.Item(32,4).Value = 0
$ws1.Cells.Item(32,5).Value = 0
$ws1.Cells.Item(32,6).Value = 2
$ws1.Cells.Item(32,7).Value = 0
$ws1.Cells.Item(32,8).Value = 3
$ws1.Cells.Item(32,9).Value = 10

$ws1.Cells.Item(33,1).Value = "Clermont"
$ws1.Cells.Item(33,2).Value = 0
$ws1.Cells.Item(33,3).Value = 0
$ws1.Cells.Item(33,4).Value = 1
$ws1.Cells.Item(33,5).Value = 0
$ws1.Cells.Item(33,6).Value = 0
$ws1.Cells.Item(33,7).Value = 2
$ws1.Cells.Item(33,8).Value = 10
$ws1.Cells.Item(33,9).Value = 3

$ws1.Cells.Item(34,1).Value = "Monaco"
$ws1.Cells.Item(34,2).Value = 1
$ws1.Cells.Item(34,3).Value = 1
$ws1.Cells.Item(34,4).Value = 0
$ws1.Cells.Item(34,5).Value = 0
$ws1.Cells.Item(34,6).Value = 3
$ws1.Cells.Item(34,7).Value = 0
$ws1.Cells.Item(34,8).Value = 4
$ws1.Cells.Item(34,9).Value = 3

$ws1.Cells.Item(35,1).Value = "Strasbourg"
$ws1.Cells.Item(35,2).Value = 0
$ws1.Cells.Item(35,3).Value = 0
$ws1.Cells.Item(35,4).Value = 1
$ws1.Cells.Item(35,5).Value = 0
$ws1.Cells.Item(35,6).Value = 0
$ws1.Cells.Item(35,7).Value = 3
$ws1.Cells.Item(35,8).Value = 3
$ws1.Cells.Item(35,9).Value = 4

$ws1.Cells.Item(36,1).Value = "Lens"
$ws1.Cells.Item(36,2).Value = 1
$ws1.Cells.Item(36,3).Value = 1
$ws1.Cells.Item(36,4).Value = 1
$ws1.Cells.Item(36,5).Value = 1
$ws1.Cells.Item(36,6).Value = 1
$ws1.Cells.Item(36,7).Value = 1
$ws1.Cells.Item(36,8).Value = 11
$ws1.Cells.Item(36,9).Value = 0

$ws1.Cells.Item(37,1).Value = "Rennes"
$ws1.Cells.Item(37,2).Value = 0
$ws1.Cells.Item(37,3).Value = 1
$ws1.Cells.Item(37,4).Value = 1
$ws1.Cells.Item(37,5).Value = 1
$ws1.Cells.Item(37,6).Value = 1
$ws1.Cells.Item(37,7).Value = 1
$ws1.Cells.Item(37,8).Value = 0
$ws1.Cells.Item(37,9).Value = 11

# --- BD_Jogo: append rows 11-19 (new round of matches) ---
$ws2.Cells.Item(11,1).Value = 1
$ws2.Cells.Item(11,2).Value = 4
$ws2.Cells.Item(11,3).Value = 13
$ws2.Cells.Item(11,4).Value = "Metz"
$ws2.Cells.Item(11,5).Value = "Marseille"

$ws2.Cells.Item(12,1).Value = 1
$ws2.Cells.Item(12,2).Value = 5
$ws2.Cells.Item(12,3).Value = 8
$ws2.Cells.Item(12,4).Value = "Lyon"
$ws2.Cells.Item(12,5).Value = "Montpelier"

$ws2.Cells.Item(13,1).Value = 1
$ws2.Cells.Item(13,2).Value = 2
$ws2.Cells.Item(13,3).Value = 10
$ws2.Cells.Item(13,4).Value = "Toulouse"
$ws2.Cells.Item(13,5).Value = "PSG"

$ws2.Cells.Item(14,1).Value = 0
$ws2.Cells.Item(14,2).Value = 2
$ws2.Cells.Item(14,3).Value = 12
$ws2.Cells.Item(14,4).Value = "Lille"
$ws2.Cells.Item(14,5).Value = "Nantes"

$ws2.Cells.Item(15,1).Value = 1
$ws2.Cells.Item(15,2).Value = 3
$ws2.Cells.Item(15,3).Value = 8
$ws2.Cells.Item(15,4).Value = "Le Havre"
$ws2.Cells.Item(15,5).Value = "Brest"

$ws2.Cells.Item(16,1).Value = 1
$ws2.Cells.Item(16,2).Value = 2
$ws2.Cells.Item(16,3).Value = 9
$ws2.Cells.Item(16,4).Value = "Lorient"
$ws2.Cells.Item(16,5).Value = "Nice"

$ws2.Cells.Item(17,1).Value = 0
$ws2.Cells.Item(17,2).Value = 2
$ws2.Cells.Item(17,3).Value = 13
$ws2.Cells.Item(17,4).Value = "Reims"
$ws2.Cells.Item(17,5).Value = "Clermont"

$ws2.Cells.Item(18,1).Value = 0
$ws2.Cells.Item(18,2).Value = 3
$ws2.Cells.Item(18,3).Value = 7
$ws2.Cells.Item(18,4).Value = "Monaco"
$ws2.Cells.Item(18,5).Value = "Strasbourg"

$ws2.Cells.Item(19,1).Value = 1
$ws2.Cells.Item(19,2).Value = 2
$ws2.Cells.Item(19,3).Value = 11
$ws2.Cells.Item(19,4).Value = "Lens"
$ws2.Cells.Item(19,5).Value = "Rennes"

